$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for refseq and crispick columns
$ws.Range("G1").Value = "refseq"
$ws.Range("H1").Value = "crispick"

$rowCount = 74
$data = New-Object 'object[,]' $rowCount,2
$data[0,0] = "NC_000011.9"
$data[0,1] = "NC_000011.9:+:47377184-47377382"
$data[1,0] = "NC_000011.9"
$data[1,1] = "NC_000011.9:+:47410789-47410987"
$data[2,0] = "NC_000011.9"
$data[2,1] = "NC_000011.9:+:47416647-47416845"
$data[3,0] = "NC_000011.9"
$data[3,1] = "NC_000011.9:+:47430500-47430698"
$data[4,0] = "NC_000011.9"
$data[4,1] = "NC_000011.9:+:59961328-59961585"
$data[5,0] = "NC_000011.9"
$data[5,1] = "NC_000011.9:+:60000475-60000673"
$data[6,0] = "NC_000011.9"
$data[6,1] = "NC_000011.9:+:60002836-60003034"
$data[7,0] = "NC_000011.9"
$data[7,1] = "NC_000011.9:+:60020013-60020211"
$data[8,0] = "NC_000011.9"
$data[8,1] = "NC_000011.9:+:60031171-60031369"
$data[9,0] = "NC_000011.9"
$data[9,1] = "NC_000011.9:+:60033272-60033547"
$data[10,0] = "NC_000011.9"
$data[10,1] = "NC_000011.9:+:85814931-85815129"
$data[11,0] = "NC_000011.9"
$data[11,1] = "NC_000011.9:+:85828449-85828650"
$data[12,0] = "NC_000011.9"
$data[12,1] = "NC_000011.9:+:85867776-85867974"
$data[13,0] = "NC_000011.9"
$data[13,1] = "NC_000011.9:+:121352978-121353176"
$data[14,0] = "NC_000008.10"
$data[14,1] = "NC_000008.10:+:27195022-27195220"
$data[15,0] = "NC_000008.10"
$data[15,1] = "NC_000008.10:+:27219888-27220086"
$data[16,0] = "NC_000008.10"
$data[16,1] = "NC_000008.10:+:27456154-27456352"
$data[17,0] = "NC_000008.10"
$data[17,1] = "NC_000008.10:+:27466058-27466414"
$data[18,0] = "NC_000008.10"
$data[18,1] = "NC_000008.10:+:145158508-145158706"
$data[19,0] = "NC_000016.9"
$data[19,1] = "NC_000016.9:+:29984740-29984938"
$data[20,0] = "NC_000016.9"
$data[20,1] = "NC_000016.9:+:30021303-30021501"
$data[21,0] = "NC_000016.9"
$data[21,1] = "NC_000016.9:+:31154047-31154245"
$data[22,0] = "NC_000006.11"
$data[22,1] = "NC_000006.11:+:32590538-32590740"
$data[23,0] = "NC_000006.11"
$data[23,1] = "NC_000006.11:+:32591489-32591687"
$data[24,0] = "NC_000006.11"
$data[24,1] = "NC_000006.11:+:114637332-114637530"
$data[25,0] = "NC_000006.11"
$data[25,1] = "NC_000006.11:+:114645138-114645336"
$data[26,0] = "NC_000017.10"
$data[26,1] = "NC_000017.10:+:1639620-1639818"
$data[27,0] = "NC_000017.10"
$data[27,1] = "NC_000017.10:+:1639866-1640133"
$data[28,0] = "NC_000017.10"
$data[28,1] = "NC_000017.10:+:1640430-1640892"
$data[29,0] = "NC_000017.10"
$data[29,1] = "NC_000017.10:+:1640936-1641134"
$data[30,0] = "NC_000017.10"
$data[30,1] = "NC_000017.10:+:1641616-1641814"
$data[31,0] = "NC_000017.10"
$data[31,1] = "NC_000017.10:+:1648403-1648601"
$data[32,0] = "NC_000017.10"
$data[32,1] = "NC_000017.10:+:18044493-18044691"
$data[33,0] = "NC_000017.10"
$data[33,1] = "NC_000017.10:+:18044999-18045197"
$data[34,0] = "NC_000017.10"
$data[34,1] = "NC_000017.10:+:18090555-18090753"
$data[35,0] = "NC_000017.10"
$data[35,1] = "NC_000017.10:+:44843037-44843235"
$data[36,0] = "NC_000017.10"
$data[36,1] = "NC_000017.10:+:44848339-44848616"
$data[37,0] = "NC_000017.10"
$data[37,1] = "NC_000017.10:+:44859616-44859814"
$data[38,0] = "NC_000017.10"
$data[38,1] = "NC_000017.10:+:44863034-44863232"
$data[39,0] = "NC_000019.9"
$data[39,1] = "NC_000019.9:+:1819135-1819333"
$data[40,0] = "NC_000019.9"
$data[40,1] = "NC_000019.9:+:49228173-49228371"
$data[41,0] = "NC_000019.9"
$data[41,1] = "NC_000019.9:+:54815267-54815465"
$data[42,0] = "NC_000007.13"
$data[42,1] = "NC_000007.13:+:54949157-54949355"
$data[43,0] = "NC_000007.13"
$data[43,1] = "NC_000007.13:+:100091696-100091894"
$data[44,0] = "NC_000007.13"
$data[44,1] = "NC_000007.13:+:143104232-143104430"
$data[45,0] = "NC_000012.11"
$data[45,1] = "NC_000012.11:+:113591338-113591536"
$data[46,0] = "NC_000012.11"
$data[46,1] = "NC_000012.11:+:113634956-113635154"
$data[47,0] = "NC_000012.11"
$data[47,1] = "NC_000012.11:+:113659522-113659851"
$data[48,0] = "NC_000012.11"
$data[48,1] = "NC_000012.11:+:113679499-113679697"
$data[49,0] = "NC_000010.10"
$data[49,1] = "NC_000010.10:+:82265172-82265370"
$data[50,0] = "NC_000010.10"
$data[50,1] = "NC_000010.10:+:82269362-82269710"
$data[51,0] = "NC_000010.10"
$data[51,1] = "NC_000010.10:+:82269749-82269947"
$data[52,0] = "NC_000010.10"
$data[52,1] = "NC_000010.10:+:98017767-98017965"
$data[53,0] = "NC_000010.10"
$data[53,1] = "NC_000010.10:+:98048166-98048364"
$data[54,0] = "NC_000010.10"
$data[54,1] = "NC_000010.10:+:124127891-124128089"
$data[55,0] = "NC_000009.11"
$data[55,1] = "NC_000009.11:+:107672266-107672464"
$data[56,0] = "NC_000020.10"
$data[56,1] = "NC_000020.10:+:54984669-54984867"
$data[57,0] = "NC_000020.10"
$data[57,1] = "NC_000020.10:+:54997469-54997667"
$data[58,0] = "NC_000020.10"
$data[58,1] = "NC_000020.10:+:55003366-55003564"
$data[59,0] = "NC_000020.10"
$data[59,1] = "NC_000020.10:+:55012893-55013091"
$data[60,0] = "NC_000020.10"
$data[60,1] = "NC_000020.10:+:55015067-55015265"
$data[61,0] = "NC_000020.10"
$data[61,1] = "NC_000020.10:+:55018057-55018359"
$data[62,0] = "NC_000014.8"
$data[62,1] = "NC_000014.8:+:53319733-53320004"
$data[63,0] = "NC_000014.8"
$data[63,1] = "NC_000014.8:+:53320429-53320627"
$data[64,0] = "NC_000014.8"
$data[64,1] = "NC_000014.8:+:53346832-53347030"
$data[65,0] = "NC_000014.8"
$data[65,1] = "NC_000014.8:+:53374678-53375066"
$data[66,0] = "NC_000002.11"
$data[66,1] = "NC_000002.11:+:127892669-127892909"
$data[67,0] = "NC_000002.11"
$data[67,1] = "NC_000002.11:+:203926172-203926370"
$data[68,0] = "NC_000005.9"
$data[68,1] = "NC_000005.9:+:86181604-86181802"
$data[69,0] = "NC_000005.9"
$data[69,1] = "NC_000005.9:+:86300359-86300557"
$data[70,0] = "NC_000005.9"
$data[70,1] = "NC_000005.9:+:156513183-156513381"
$data[71,0] = "NC_000015.9"
$data[71,1] = "NC_000015.9:+:64382798-64383128"
$data[72,0] = "NC_000021.8"
$data[72,1] = "NC_000021.8:+:27541845-27542043"
$data[73,0] = "NC_000004.11"
$data[73,1] = "NC_000004.11:+:40198747-40198945"

$rng = $ws.Range("G2:H" + (1 + $rowCount))
$rng.Value = $data

Write-Host "Done writing refseq/crispick columns."